# Auto-generated script to apply scheduled-runner value updates to Sheets/Ixion_Profits.xlsx
# Each block updates the market/profit columns (H:N) for a specific leve row on a specific sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1827.225
$ws.Range("I70").Value = 1929.2
$ws.Range("J70").Value = 1725.25
$ws.Range("K70").Value = 5787.6
$ws.Range("L70").Value = 5175.75
$ws.Range("M70").Value = -5517.6
$ws.Range("N70").Value = -5715.75

$ws.Range("H73").Value = 1827.225
$ws.Range("I73").Value = 1929.2
$ws.Range("J73").Value = 1725.25
$ws.Range("K73").Value = 5787.6
$ws.Range("L73").Value = 5175.75
$ws.Range("M73").Value = -4851.6
$ws.Range("N73").Value = -7047.75

$ws.Range("H76").Value = 3727
$ws.Range("I76").Value = 3310.4666
$ws.Range("J76").Value = 5289
$ws.Range("K76").Value = 3310.4666
$ws.Range("L76").Value = 5289
$ws.Range("M76").Value = -2995.4666
$ws.Range("N76").Value = -5919.1113

$ws.Range("H79").Value = 3727
$ws.Range("I79").Value = 3310.4666
$ws.Range("J79").Value = 5289
$ws.Range("K79").Value = 3310.4666
$ws.Range("L79").Value = 5289
$ws.Range("M79").Value = -2218.4666
$ws.Range("N79").Value = -7473.1113

$ws.Range("H138").Value = 2973.1016
$ws.Range("I138").Value = 1315.9062
$ws.Range("J138").Value = 4406.3516
$ws.Range("K138").Value = 3947.7186
$ws.Range("L138").Value = 13219.0548
$ws.Range("M138").Value = 1192.2814
$ws.Range("N138").Value = -23499.0548

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H141").Value = 2112.15
$ws.Range("I141").Value = 2112.15
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6336.450000000001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1156.450000000001
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 100002260
$ws.Range("I63").Value = 100002260
$ws.Range("K63").Value = 100002260
$ws.Range("M63").Value = -100001574

$ws.Range("H66").Value = 100002260
$ws.Range("I66").Value = 100002260
$ws.Range("K66").Value = 500011300
$ws.Range("M66").Value = -500007868

$ws.Range("H68").Value = 35050
$ws.Range("I68").Value = 28000
$ws.Range("J68").Value = 42100
$ws.Range("K68").Value = 28000
$ws.Range("L68").Value = 42100
$ws.Range("M68").Value = -27189
$ws.Range("N68").Value = -43722

$ws.Range("H71").Value = 35050
$ws.Range("I71").Value = 28000
$ws.Range("J71").Value = 42100
$ws.Range("K71").Value = 84000
$ws.Range("L71").Value = 126300
$ws.Range("M71").Value = -79944
$ws.Range("N71").Value = -134412

$ws.Range("H122").Value = 952086.25
$ws.Range("I122").Value = 1070846.5
$ws.Range("K122").Value = 3212539.5
$ws.Range("M122").Value = -3210089.5

$ws.Range("H132").Value = 2674.225
$ws.Range("I132").Value = 1427.8064
$ws.Range("J132").Value = 6967.4443
$ws.Range("K132").Value = 4283.4192
$ws.Range("L132").Value = 20902.3329
$ws.Range("M132").Value = -1753.4192
$ws.Range("N132").Value = -25962.3329

$ws.Range("H139").Value = 40426.75
$ws.Range("J139").Value = 40426.75
$ws.Range("L139").Value = 40426.75
$ws.Range("N139").Value = -50706.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 22788.924
$ws.Range("I82").Value = 7564.25
$ws.Range("J82").Value = 29555.445
$ws.Range("K82").Value = 7564.25
$ws.Range("L82").Value = 29555.445
$ws.Range("M82").Value = -7181.25
$ws.Range("N82").Value = -30321.445

$ws.Range("H85").Value = 22788.924
$ws.Range("I85").Value = 7564.25
$ws.Range("J85").Value = 29555.445
$ws.Range("K85").Value = 7564.25
$ws.Range("L85").Value = 29555.445
$ws.Range("M85").Value = -6238.25
$ws.Range("N85").Value = -32207.445

$ws.Range("H138").Value = 59795
$ws.Range("J138").Value = 59795
$ws.Range("L138").Value = 59795
$ws.Range("N138").Value = -70075

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9632398
$ws.Range("I99").Value = 18842.4
$ws.Range("J99").Value = 15640870
$ws.Range("K99").Value = 18842.4
$ws.Range("L99").Value = 15640870
$ws.Range("M99").Value = -17344.4
$ws.Range("N99").Value = -15643866

$ws.Range("H126").Value = 9632398
$ws.Range("I126").Value = 18842.4
$ws.Range("J126").Value = 15640870
$ws.Range("K126").Value = 56527.2
$ws.Range("L126").Value = 46922610
$ws.Range("M126").Value = -54057.2
$ws.Range("N126").Value = -46927550

$ws.Range("H138").Value = 34814.75
$ws.Range("J138").Value = 34814.75
$ws.Range("L138").Value = 34814.75
$ws.Range("N138").Value = -45094.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 3833.3333
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 3833.3333
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 11499.9999
$ws.Range("N76").Value = -12265.9999
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 3833.3333
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 3833.3333
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 11499.9999
$ws.Range("N79").Value = -14151.9999
$ws.Range("M79").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6180.8423
$ws.Range("I70").Value = 6297.5713
$ws.Range("K70").Value = 6297.5713
$ws.Range("M70").Value = -6027.5713

$ws.Range("H73").Value = 6180.8423
$ws.Range("I73").Value = 6297.5713
$ws.Range("K73").Value = 6297.5713
$ws.Range("M73").Value = -5361.5713

$ws.Range("H80").Value = 3800
$ws.Range("J80").Value = 3666.6667
$ws.Range("L80").Value = 3666.6667
$ws.Range("N80").Value = -5662.6667

$ws.Range("H83").Value = 3800
$ws.Range("J83").Value = 3666.6667
$ws.Range("L83").Value = 18333.3335
$ws.Range("N83").Value = -28317.3335

$ws.Range("H102").Value = 2726
$ws.Range("I102").Value = 978.2857
$ws.Range("K102").Value = 978.2857
$ws.Range("M102").Value = 643.7143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4632113
$ws.Range("I22").Value = 27778954
$ws.Range("J22").Value = 2745
$ws.Range("K22").Value = 27778954
$ws.Range("L22").Value = 2745
$ws.Range("M22").Value = -27778659
$ws.Range("N22").Value = -3335

$ws.Range("H27").Value = 4632113
$ws.Range("I27").Value = 27778954
$ws.Range("J27").Value = 2745
$ws.Range("K27").Value = 27778954
$ws.Range("L27").Value = 2745
$ws.Range("M27").Value = -27778847
$ws.Range("N27").Value = -2959

$ws.Range("H40").Value = 62502716
$ws.Range("I40").Value = 76925560
$ws.Range("K40").Value = 76925560
$ws.Range("M40").Value = -76925424

$ws.Range("H63").Value = 17750
$ws.Range("I63").Value = 13000
$ws.Range("K63").Value = 13000
$ws.Range("M63").Value = -12251

$ws.Range("H66").Value = 17750
$ws.Range("I66").Value = 13000
$ws.Range("K66").Value = 39000
$ws.Range("M66").Value = -35256

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 7559.8335
$ws.Range("I45").Value = 3069
$ws.Range("J45").Value = 9805.25
$ws.Range("K45").Value = 3069
$ws.Range("L45").Value = 9805.25
$ws.Range("M45").Value = -2578
$ws.Range("N45").Value = -10787.25

$ws.Range("H70").Value = 32000
$ws.Range("J70").Value = 32000
$ws.Range("L70").Value = 32000
$ws.Range("N70").Value = -32630

$ws.Range("H73").Value = 32000
$ws.Range("J73").Value = 32000
$ws.Range("L73").Value = 32000
$ws.Range("N73").Value = -34184

$ws.Range("H93").Value = 23250
$ws.Range("J93").Value = 23250
$ws.Range("L93").Value = 23250
$ws.Range("N93").Value = -28242

$ws.Range("H106").Value = 34000
$ws.Range("J106").Value = 34000
$ws.Range("L106").Value = 34000
$ws.Range("N106").Value = -36524

$ws.Range("H122").Value = 3131.0286
$ws.Range("I122").Value = 2797.2
$ws.Range("J122").Value = 5134
$ws.Range("K122").Value = 8391.599999999999
$ws.Range("L122").Value = 15402
$ws.Range("M122").Value = -5941.599999999999
$ws.Range("N122").Value = -20302
